$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "XMOS Dualchip Planning": build the By-Core (X0..X3) 1-bit
# port usage table mirrored alongside the existing By-Block table.
# ------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("XMOS Dualchip Planning")

# New section labels first (matches authoring order of the new strings).
$ws6.Range("A16").Value = "By Block ^"
$ws6.Range("G16").Value = "By Core^"

$ws6.Range("G1").ColumnWidth = 18.5
$ws6.Range("G1").Value = "Block"
$ws6.Range("H1").Value = "X0 1-bit"
$ws6.Range("I1").Value = "X1 1-bit"
$ws6.Range("J1").Value = "X2 1-bit"
$ws6.Range("K1").Value = "X3 1-bit"

$ws6.Range("G2").Value = "Pmod 0"
$ws6.Range("G3").Value = "Pmod 1"
$ws6.Range("G4").Value = "Pmod 2"
$ws6.Range("G5").Value = "Gadgeteer"
$ws6.Range("G6").Value = "SPI"
$ws6.Range("H6").Value = 5
$ws6.Range("G7").Value = "I2C"
$ws6.Range("G8").Value = "PWM"
$ws6.Range("J8").Value = 0
$ws6.Range("G9").Value = "Enc"
$ws6.Range("G10").Value = "STM32 UART"
$ws6.Range("G11").Value = "Xbee UART"

# PWM block's chip assignment moved from U1 to U2; drop the stray "2?" note.
$ws6.Range("D8").Value = "U2"
$ws6.Range("F8").ClearContents()

$ws6.Range("G13").Value = "Total 1-bits"
$ws6.Range("H13").Formula = "=SUM(H3:H11)"
$ws6.Range("I13:K13").Formula = "=SUM(I2:I11)"

$ws6.Range("I5").Select()

# ------------------------------------------------------------------
# Sheet "XMOS Dualchip": assign the microSD SS line to X0P1K
# (moves the uSD_SS label from the stray E37 cell to D36, the X0
# Signal column for P1K).
# ------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("XMOS Dualchip")

$ws5.Range("D36").Value = "uSD_SS"
$ws5.Range("E37").ClearContents()

$ws5.Activate()
$ws5.Range("D36").Select()
